# Page items put under different pages
#
# The "county"/"India" field (column D) on the LoginTest page is replaced
# by a "country"/"Italy" field. Delete the old column then insert a fresh
# one in its place (rather than just overwriting the values) so the new
# cells pick up the neighbouring column's formatting the same way Excel's
# own insert-column flow does (this is what gives D2 the same style as the
# hyperlinked C2 cell).
$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("LoginTest")
$wsLogin.Columns.Item(4).Delete()
$wsLogin.Columns.Item(4).Insert()
$wsLogin.Range("D1").Value = "country"
$wsLogin.Range("D2").Value = "Italy"

# Column C (email) also ended up a bit wider while this was being edited.
$wsLogin.Columns.Item(3).ColumnWidth = 26.83

# LoginTest becomes the active/selected page (it was NewCarsTest before).
$wsLogin.Activate()
[void]$wsLogin.Range("D7").Select()
